# "Last Day at R Systems" - refresh the PO_Detail sheet with the current
# batch of POs stuck on "UNTPRG-Confirm PO".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PO_Detail")

# The old single-PO sample (row 2) is replaced by a small table: row 2 keeps
# its "1" marker/PO/note, rows 3-6 are new "2" marked POs.
$ws.Range("A2").ClearContents()

$ws.Range("B2").Value = "'1"
$ws.Range("C2").Value = "'01000996"
$ws.Range("D2").Value = "'Stuck on UNTPRG-Confirm PO"

$ws.Range("B3").Value = "'2"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "'01000949"

$ws.Range("B4").Value = "'2"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "'01000942"

$ws.Range("B5").Value = "'2"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "'01000943"

$ws.Range("B6").Value = "'2"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "'01001011"
$ws.Range("D6").Value = "'Stuck on UNTPRG-Confirm PO"

$ws.Range("C3").Select() | Out-Null
